$wb = $excel.ActiveWorkbook

# --- Update timestamps on existing "data" sheet ---
$ws1 = $wb.Worksheets.Item("data")
$ws1.Range("F2").Value = "2021-10-05 14:19:39.248336"
$ws1.Range("F3").Value = "2021-10-05 14:19:39.248343"

# --- Add new "metadata" sheet after "data" ---
$newWs = $wb.Worksheets.Add($null, $ws1)
$newWs.Name = "metadata"

# Header row (row 1), columns B:G
$newWs.Range("B1").Value = "data_name"
$newWs.Range("C1").Value = "data_id"
$newWs.Range("D1").Value = "data_version"
$newWs.Range("E1").Value = "data_version_created"
$newWs.Range("F1").Value = "panel_query_time"
$newWs.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "Combined factor V and VIII deficiency"
$newWs.Range("C2").Value = 517
$newWs.Range("D2").NumberFormat = "@"
$newWs.Range("D2").Value = "1.6"
$newWs.Range("D2").ClearFormats()
$newWs.Range("E2").Value = "2021-03-23T09:59:04.457039Z"
$newWs.Range("F2").Value = "2021-10-05 14:19:39.244619"
$newWs.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/517/?format=json"

# Match the header/index styling used on the "data" sheet (bold, centered, bordered)
$ws1.Range("B1").Copy()
$newWs.Range("B1:G1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$newWs.Range("A2").PasteSpecial(-4122)

# Restore original active sheet/view state
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
